$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.531.00"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "3.938.38"
$ws.Range("E3").Value = "  +3.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "473.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000347"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").Value = "4.572.16"
$ws.Range("E13").Value = "  +4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "3.952.94"
$ws.Range("E16").Value = "  +6.19%  "

$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "67.746.13"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "719.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.38%  "

$ws.Range("E32").Value = "  +2.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.92%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.51%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0826"
$ws.Range("E35").Value = "  +20.97%  "

$ws.Range("E36").Value = "  -3.77%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("E40").Value = "  +4.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +2.06%  "

$ws.Range("E44").Value = "  +6.23%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.26%  "

$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "147.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.49%  "

$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.86%  "
